# Leakage test / Jankowski comparison samples
# Finished processing Nov 2014 flux info (O'Connell vs Jankowski protocols):
# fill in the previously-placeholder "S3" chamber rows (26-45) with their
# real sample date and sample names, and clear the highlight used to mark
# them as still-to-do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sampleDate = "2014.11.14"
$chambers = @("A", "B", "C", "D", "E")
$timePts = @("T0", "T15", "T30", "T45")

$row = 26
foreach ($chamber in $chambers) {
    foreach ($timePt in $timePts) {
        $ws.Range("D$row").Value = $sampleDate
        $ws.Range("E$row").Value = "S3-$chamber-14-$timePt"

        # Clear the "still needs data" highlight fill on both cells.
        $ws.Range("D$row").Interior.Pattern = -4142
        $ws.Range("D$row").Interior.ColorIndex = -4142
        $ws.Range("E$row").Interior.Pattern = -4142
        $ws.Range("E$row").Interior.ColorIndex = -4142

        $row = $row + 1
    }
}

# Update view state to reflect where the editor was last working.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E30").Select()
